# Updated symbol list on Sun Feb 12 15:39:05 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns on Sheet1 with new
# quote data. Source values are stored as literal text (e.g. "310.21",
# "0.38%") rather than numbers, so each assignment uses a leading
# apostrophe to force Excel to keep the exact text instead of re-parsing
# it as a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.21"
$ws.Range("E2").Value = "'0.38%"

$ws.Range("D3").Value = "'41.10"
$ws.Range("E3").Value = "'-0.43%"

$ws.Range("D4").Value = "'5.195"
$ws.Range("E4").Value = "'1.39%"

$ws.Range("D5").Value = "'0.07684"
$ws.Range("E5").Value = "'0.43%"

$ws.Range("D6").Value = "'4.296"
$ws.Range("E6").Value = "'0.49%"

$ws.Range("D7").Value = "'1.696"
$ws.Range("E7").Value = "'4.89%"

$ws.Range("D8").Value = "'0.9406"
$ws.Range("E8").Value = "'3.47%"

$ws.Range("D9").Value = "'2.423"
$ws.Range("E9").Value = "'-2.02%"

$ws.Range("D10").Value = "'0.1279"
$ws.Range("E10").Value = "'7.71%"

$ws.Range("D11").Value = "'0.1835"
$ws.Range("E11").Value = "'1.79%"

$ws.Range("D12").Value = "'0.09236"
$ws.Range("E12").Value = "'0.73%"

$ws.Range("D13").Value = "'0.04242"
$ws.Range("E13").Value = "'-0.22%"

$ws.Range("D14").Value = "'0.1052"

$ws.Range("D15").Value = "'0.001282"
$ws.Range("E15").Value = "'2.56%"

$ws.Range("D16").Value = "'0.005890"
$ws.Range("E16").Value = "'0.27%"

$ws.Range("D19").Value = "'7.510"
$ws.Range("E19").Value = "'8.55%"

$ws.Range("E20").Value = "'-1.83%"

$ws.Range("D21").Value = "'0.2721"
$ws.Range("E21").Value = "'-0.54%"

$ws.Range("D22").Value = "'0.04011"
$ws.Range("E22").Value = "'-1.18%"

$ws.Range("E23").Value = "'-0.35%"

$ws.Range("D24").Value = "'0.004247"
$ws.Range("E24").Value = "'5.59%"

$ws.Range("E25").Value = "'0.08%"

$ws.Range("E38").Value = "'4.30%"

$ws.Range("D39").Value = "'0.05318"
$ws.Range("E39").Value = "'1.28%"

$ws.Range("D40").Value = "'0.007846"
$ws.Range("E40").Value = "'0.87%"

$ws.Range("E41").Value = "'0.88%"

$ws.Range("D42").Value = "'0.006662"
$ws.Range("E42").Value = "'-1.74%"

$ws.Range("D43").Value = "'0.001941"
$ws.Range("E43").Value = "'-0.38%"

$ws.Range("D44").Value = "'0.007425"
$ws.Range("E44").Value = "'-1.64%"

$ws.Range("E45").Value = "'0.49%"

$ws.Range("D46").Value = "'0.00006770"
$ws.Range("E46").Value = "'-1.84%"

$ws.Range("E47").Value = "'0.08%"

$ws.Range("D48").Value = "'0.2179"
$ws.Range("E48").Value = "'179.62%"

$ws.Range("E49").Value = "'3.46%"

$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.08%"

$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.08%"
